# "update round 2 move list to account for movement"
#
# On the "Template" sheet, the move list repeats a 10-row cycle (rows 4..300).
# Within each cycle the 4th row (column B) and the 9th row (column C) used to
# always say "Wait". They are updated so that movement is now represented
# explicitly, alternating between "Move:down" and "Move:up" every cycle so the
# party moves down a level then back up on the following loop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$moveNames = @("Move:down", "Move:up")

# Column B (index 2) on rows 4, 14, 24, ... 294
$i = 0
for ($row = 4; $row -le 300; $row += 10) {
    $ws.Cells.Item($row, 2).Value = $moveNames[$i % 2]
    $i++
}

# Column C (index 3) on rows 9, 19, 29, ... 299 (opposite phase of column B)
$i = 0
for ($row = 9; $row -le 300; $row += 10) {
    $ws.Cells.Item($row, 3).Value = $moveNames[1 - ($i % 2)]
    $i++
}
